# C5-PowerPoint.pptx edit — Mon, Jul 27, 2020  2:05:56 PM
#
# 1) Slide 6's table ("Google Shape;127;p18") is switched to a different
#    built-in PowerPoint table style (Table Design gallery).
# 2) The deck's design/theme colour palette is swapped from the "Integral"
#    palette to the standard "Office" palette (done by rewriting each of
#    the 12 theme colour slots through ThemeColorScheme, which is the
#    supported way to edit the active theme's colours from the object
#    model).

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $null
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $candidate = $tableSlide.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
    }
}
$tableShape.Table.ApplyStyle("{A6E20BB9-44FC-4921-88AE-04C7E6C82557}")

# --- 2) Theme colours: Integral -> Office ---------------------------------
function HexToBgrInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Colors(1..12) order == dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$colorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $colorScheme.Colors($i).RGB = HexToBgrInt $officeThemeColors[$i - 1]
}
